$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Статус" column (E) is being moved to the front of the table (new column A),
# and is given a background color depending on its value ("Нет" -> green, "Да" -> red).
# "Ссылка" (F) stays as the last column.

# 1) Insert a new blank column before column A; this shifts A:F -> B:G
#    (Company/INN/Start/End/Status/Link)
$ws.Columns("A:A").Insert()

# 2) Copy the Status column (now F, after the shift) into the new column A
$ws.Range("A1:A5").Value2 = $ws.Range("F1:F5").Value2

# 3) Remove the now-duplicated Status column (F); this shifts the trailing
#    Link column (G) back into F
$ws.Columns("F:F").Delete()

# 4) Color the moved Status cells based on their value: green for "Нет", red for "Да"
$ws.Range("A2").Interior.Color = 32768
$ws.Range("A3").Interior.Color = 32768
$ws.Range("A4").Interior.Color = 255
$ws.Range("A5").Interior.Color = 255
